# Add a new "ToDo" list item ("Stop skipping player 1's first turn") right
# after the "...LifeDetails' Excel sheet" item, keeping the same
# ListParagraph style / numbering (numId 4) used by the other items, and
# move the trailing "_GoBack" bookmark along so it still ends up at the
# very end of the (now last) list item.

$d = $word.ActiveDocument

$newItemText = "Stop skipping player 1" + [string][char]0x2019 + "s first turn"

# Locate the end of the "...Excel sheet" list item (the paragraph that
# currently carries the _GoBack bookmark right after its text).
$findRange = $d.Content
$found = $findRange.Find.Execute("Excel sheet", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Excel sheet' to-do item to append after."
}
$findRange.Collapse(0)
$insertPos = $findRange.Start

# Insert a new paragraph right after that item. It inherits the
# ListParagraph style + numbering automatically since it is split off of
# that list paragraph.
$splitRange = $d.Range($insertPos, $insertPos)
$splitRange.InsertParagraphAfter()

# The new paragraph now starts exactly one position after $insertPos
# (the paragraph mark that was just inserted there).
$newParaStart = $insertPos + 1
$insAt = $d.Range($newParaStart, $newParaStart)
$insAt.InsertAfter($newItemText)

# --- Move the _GoBack bookmark to the end of the new last paragraph ---
# Work out where the end of the freshly-typed text now is.
$newTextEnd = $newParaStart + $newItemText.Length

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()

    # Adding a bookmark whose position coincides with the very end of the
    # document's content is unreliable in this runtime, so pad the
    # document with one throwaway character first, add the bookmark at
    # its real (now no-longer-final) position, then remove the padding.
    $tailRange = $d.Content
    $tailRange.Collapse(0)
    $tailRange.InsertAfter("X")

    $bmRange = $d.Range($newTextEnd, $newTextEnd)
    $d.Bookmarks.Add("_GoBack", $bmRange)

    $padRange = $d.Range($d.Content.End - 2, $d.Content.End - 1)
    $padRange.Delete()
}
